{"js": "// Update the two-digit multiplication problems in the document.\n// Each cell contains a single \"AA\u00d7BB=\" expression; replace the operands\n// per the mapping below using Word's search API (exact, case-sensitive,\n// no wildcards) and Range.insertText(..., \"Replace\").\n\nconst replacements = [\n  [\"71\u00d719=\", \"79\u00d759=\"],\n  [\"75\u00d774=\", \"25\u00d791=\"],\n  [\"69\u00d752=\", \"88\u00d793=\"],\n  [\"92\u00d740=\", \"49\u00d764=\"],\n  [\"35\u00d724=\", \"86\u00d715=\"],\n  [\"92\u00d767=\", \"26\u00d758=\"],\n  [\"97\u00d731=\", \"30\u00d728=\"],\n  [\"96\u00d784=\", \"95\u00d763=\"],\n  [\"89\u00d720=\", \"40\u00d724=\"],\n  [\"50\u00d754=\", \"45\u00d798=\"],\n  [\"67\u00d712=\", \"57\u00d732=\"],\n  [\"83\u00d785=\", \"79\u00d734=\"],\n  [\"22\u00d717=\", \"33\u00d779=\"],\n  [\"25\u00d758=\", \"52\u00d725=\"],\n  [\"68\u00d756=\", \"30\u00d759=\"],\n  [\"64\u00d739=\", \"59\u00d774=\"],\n  [\"30\u00d776=\", \"65\u00d776=\"],\n  [\"15\u00d738=\", \"83\u00d764=\"],\n  [\"25\u00d771=\", \"67\u00d715=\"],\n  [\"47\u00d784=\", \"55\u00d728=\"],\n  [\"44\u00d771=\", \"36\u00d787=\"],\n  [\"35\u00d788=\", \"93\u00d738=\"],\n  [\"24\u00d761=\", \"13\u00d716=\"],\n  [\"30\u00d754=\", \"20\u00d739=\"],\n  [\"67\u00d732=\", \"72\u00d718=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication problems in the document.\n# Each table cell contains a single \"AA\u00d7BB=\" expression; replace the\n# operands per the mapping below using Find/Execute (exact match, no\n# wildcards) against the whole document content.\n\n$d = $word.ActiveDocument\n\n$olds = @(\"71\u00d719=\", \"75\u00d774=\", \"69\u00d752=\", \"92\u00d740=\", \"35\u00d724=\", \"92\u00d767=\", \"97\u00d731=\", \"96\u00d784=\", \"89\u00d720=\", \"50\u00d754=\", \"67\u00d712=\", \"83\u00d785=\", \"22\u00d717=\", \"25\u00d758=\", \"68\u00d756=\", \"64\u00d739=\", \"30\u00d776=\", \"15\u00d738=\", \"25\u00d771=\", \"47\u00d784=\", \"44\u00d771=\", \"35\u00d788=\", \"24\u00d761=\", \"30\u00d754=\", \"67\u00d732=\")\n$news = @(\"79\u00d759=\", \"25\u00d791=\", \"88\u00d793=\", \"49\u00d764=\", \"86\u00d715=\", \"26\u00d758=\", \"30\u00d728=\", \"95\u00d763=\", \"40\u00d724=\", \"45\u00d798=\", \"57\u00d732=\", \"79\u00d734=\", \"33\u00d779=\", \"52\u00d725=\", \"30\u00d759=\", \"59\u00d774=\", \"65\u00d776=\", \"83\u00d764=\", \"67\u00d715=\", \"55\u00d728=\", \"36\u00d787=\", \"93\u00d738=\", \"13\u00d716=\", \"20\u00d739=\", \"72\u00d718=\")\n\nfor ($i = 0; $i -lt $olds.Count; $i++) {\n    $old = $olds[$i]\n    $new = $news[$i]\n\n    $find = $d.Content.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
